$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Qty executed upto date (column C) - numeric values
$ws.Range("C8").Value = 83
$ws.Range("C9").Value = 44
$ws.Range("C10").Value = 17
$ws.Range("C11").Value = 50
$ws.Range("C12").Value = 13
$ws.Range("C13").Value = 4
$ws.Range("C14").Value = 16
$ws.Range("C15").Value = 59
$ws.Range("C16").Value = 54
$ws.Range("C17").Value = 12

# Upto date Amount (column G) - stored as text-like strings (e.g. "11264.00")
# Format as Text first so the numeric-looking string isn't auto-converted to a number.
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "11264.00"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "8024.00"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "33100.00"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "544.00"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "368.00"

# Grand Total rows
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "53300.00"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "53300.00"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "53300.00"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "53300.00"
